$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text (string) cells that look numeric are not auto-converted to numbers,
# matching the original workbook where every data cell is stored as an explicit string.

$ws.Range("D2").Value = "64.061.42"
$ws.Range("E2").Value = "  -3.39%  "

$ws.Range("D3").Value = "3.144.24"
$ws.Range("E3").Value = "  -8.89%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.09"
$ws.Range("E5").Value = "  -3.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.16"
$ws.Range("E6").Value = "  -5.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -1.29%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "3.139.80"
$ws.Range("E9").Value = "  -8.94%  "

$ws.Range("E10").Value = "  -7.12%  "

$ws.Range("E11").Value = "  -6.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  -6.24%  "

$ws.Range("D13").Value = "3.691.83"
$ws.Range("E13").Value = "  -8.75%  "

$ws.Range("E14").Value = "  +0.74%  "

$ws.Range("E15").Value = "  -9.96%  "

$ws.Range("D16").Value = "64.026.23"
$ws.Range("E16").Value = "  -3.35%  "

$ws.Range("E17").Value = "  -6.41%  "

$ws.Range("D18").Value = "3.146.31"
$ws.Range("E18").Value = "  -8.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.68"
$ws.Range("E19").Value = "  -4.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.80"
$ws.Range("E20").Value = "  -7.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.72"
$ws.Range("E21").Value = "  -5.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.12"
$ws.Range("E22").Value = "  -6.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.86"
$ws.Range("E24").Value = "  -7.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.497"
$ws.Range("E25").Value = "  -7.13%  "

$ws.Range("E26").Value = "  -8.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.51"
$ws.Range("E27").Value = "  -4.33%  "

$ws.Range("E28").Value = "  -1.51%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.88"
$ws.Range("E31").Value = "  -6.09%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("E32").Value = "  -8.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.77"
$ws.Range("E33").Value = "  -7.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.55"
$ws.Range("E34").Value = "  -7.17%  "

$ws.Range("E35").Value = "  -6.64%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "153.73"
$ws.Range("E36").Value = "  -4.95%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.41"
$ws.Range("E37").Value = "  -9.15%  "

$ws.Range("E38").Value = "  -7.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.97"
$ws.Range("E39").Value = "  -6.81%  "

$ws.Range("E40").Value = "  -6.76%  "

$ws.Range("E41").Value = "  -4.48%  "

$ws.Range("D42").Value = "2.591.38"
$ws.Range("E42").Value = "  -6.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.14"
$ws.Range("E43").Value = "  -8.09%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.25"
$ws.Range("E44").Value = "  -1.72%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.94"
$ws.Range("E45").Value = "  -7.70%  "

$ws.Range("E46").Value = "  -7.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.50"
$ws.Range("E47").Value = "  -6.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "312.87"
$ws.Range("E48").Value = "  -7.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0269"
$ws.Range("E49").Value = "  -6.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  -3.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.07%  "
